# Daily attendance processing - 2026-01-17 15:35:08
#
# Normalizes the "Recorded By" column (G) of the session-analysis sheet:
#   - "system, System, backup@backdoor.com" -> "System, system, backup@backdoor.com"
#   - "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   - "backup@backdoor.com, System"         -> "System, backup@backdoor.com"
# All other "Recorded By" values (single author, different author lists, etc.)
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$recordedByCol = 7  # column G

$lastRow = $ws.Cells.Item($ws.Rows.Count, $recordedByCol).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $current = $cell.Text

    if ($current -eq "system, System, backup@backdoor.com") {
        $cell.Value = "System, system, backup@backdoor.com"
    }
    elseif ($current -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($current -eq "backup@backdoor.com, System") {
        $cell.Value = "System, backup@backdoor.com"
    }
}
